$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize existing name in A2 to lowercase
$ws.Range("A2").Value = "sarry eldeen mohamed faisal"

# Add new name in A3
$ws.Range("A3").Value = "Mariam hany adli Hassan"

# Update selection to match the new state (A3)
$ws.Range("A3").Select()
